# SemesterGradesSheet Template refactor:
# Append the semester-grade text placeholder alongside the numeric grade
# placeholder in cell F11, and leave the selection on that cell (matching
# the author's final cursor position after editing it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F11").Value = "   {SemesterGrade} ({SemesterGradeText})"

$ws.Range("F11").Select() | Out-Null
